$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-11-01 Friday", $true, $true, $false, $false, $false, $true, 1, $false, "2024-11-02 Saturday", 2) | Out-Null
$d.Content.Find.Execute("32×64=2048", $true, $true, $false, $false, $false, $true, 1, $false, "73×16=1168", 2) | Out-Null
$d.Content.Find.Execute("58×91=5278", $true, $true, $false, $false, $false, $true, 1, $false, "17×61=1037", 2) | Out-Null
$d.Content.Find.Execute("27×57=1539", $true, $true, $false, $false, $false, $true, 1, $false, "90×29=2610", 2) | Out-Null
$d.Content.Find.Execute("47×48=2256", $true, $true, $false, $false, $false, $true, 1, $false, "24×83=1992", 2) | Out-Null
$d.Content.Find.Execute("40×17=680", $true, $true, $false, $false, $false, $true, 1, $false, "29×83=2407", 2) | Out-Null
$d.Content.Find.Execute("37×83=3071", $true, $true, $false, $false, $false, $true, 1, $false, "37×41=1517", 2) | Out-Null
$d.Content.Find.Execute("92×80=7360", $true, $true, $false, $false, $false, $true, 1, $false, "57×52=2964", 2) | Out-Null
$d.Content.Find.Execute("50×24=1200", $true, $true, $false, $false, $false, $true, 1, $false, "76×68=5168", 2) | Out-Null
$d.Content.Find.Execute("79×66=5214", $true, $true, $false, $false, $false, $true, 1, $false, "46×79=3634", 2) | Out-Null
$d.Content.Find.Execute("51×76=3876", $true, $true, $false, $false, $false, $true, 1, $false, "52×54=2808", 2) | Out-Null
$d.Content.Find.Execute("13×30=390", $true, $true, $false, $false, $false, $true, 1, $false, "87×63=5481", 2) | Out-Null
$d.Content.Find.Execute("30×72=2160", $true, $true, $false, $false, $false, $true, 1, $false, "83×60=4980", 2) | Out-Null
$d.Content.Find.Execute("31×73=2263", $true, $true, $false, $false, $false, $true, 1, $false, "34×87=2958", 2) | Out-Null
$d.Content.Find.Execute("96×77=7392", $true, $true, $false, $false, $false, $true, 1, $false, "93×41=3813", 2) | Out-Null
$d.Content.Find.Execute("66×55=3630", $true, $true, $false, $false, $false, $true, 1, $false, "24×56=1344", 2) | Out-Null
$d.Content.Find.Execute("75×35=2625", $true, $true, $false, $false, $false, $true, 1, $false, "29×83=2407", 2) | Out-Null
$d.Content.Find.Execute("77×20=1540", $true, $true, $false, $false, $false, $true, 1, $false, "38×96=3648", 2) | Out-Null
$d.Content.Find.Execute("62×84=5208", $true, $true, $false, $false, $false, $true, 1, $false, "89×27=2403", 2) | Out-Null
$d.Content.Find.Execute("93×66=6138", $true, $true, $false, $false, $false, $true, 1, $false, "63×59=3717", 2) | Out-Null
$d.Content.Find.Execute("14×47=658", $true, $true, $false, $false, $false, $true, 1, $false, "55×19=1045", 2) | Out-Null
$d.Content.Find.Execute("66×51=3366", $true, $true, $false, $false, $false, $true, 1, $false, "81×98=7938", 2) | Out-Null
$d.Content.Find.Execute("74×92=6808", $true, $true, $false, $false, $false, $true, 1, $false, "87×77=6699", 2) | Out-Null
$d.Content.Find.Execute("48×78=3744", $true, $true, $false, $false, $false, $true, 1, $false, "81×19=1539", 2) | Out-Null
$d.Content.Find.Execute("91×21=1911", $true, $true, $false, $false, $false, $true, 1, $false, "32×93=2976", 2) | Out-Null
$d.Content.Find.Execute("78×97=7566", $true, $true, $false, $false, $false, $true, 1, $false, "40×26=1040", 2) | Out-Null

$d.Save()
